$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 499.5
$ws.Range("J9").Value = 199
$ws.Range("L9").Value = 199
$ws.Range("N9").Value = -537

$ws.Range("H116").Value = 6919.8438
$ws.Range("I116").Value = 8117.591
$ws.Range("J116").Value = 4284.8
$ws.Range("K116").Value = 8117.591
$ws.Range("L116").Value = 4284.8
$ws.Range("M116").Value = -4675.591
$ws.Range("N116").Value = -11168.8

$ws.Range("H132").Value = 3343.2058
$ws.Range("I132").Value = 1841.5161
$ws.Range("J132").Value = 18860.666
$ws.Range("K132").Value = 5524.5483
$ws.Range("L132").Value = 56581.99800000001
$ws.Range("M132").Value = -2994.5483
$ws.Range("N132").Value = -61641.99800000001

$ws.Range("H138").Value = 2068.9167
$ws.Range("I138").Value = 1534.909
$ws.Range("J138").Value = 2520.7693
$ws.Range("K138").Value = 4604.727000000001
$ws.Range("L138").Value = 7562.3079
$ws.Range("M138").Value = 535.2729999999992
$ws.Range("N138").Value = -17842.3079

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2562.3125
$ws.Range("I32").Value = 2562.3125
$ws.Range("K32").Value = 2562.3125
$ws.Range("M32").Value = -2275.3125

$ws.Range("H61").Value = 2217.0715
$ws.Range("I61").Value = 2019.5454
$ws.Range("K61").Value = 2019.5454
$ws.Range("M61").Value = -1807.5454

$ws.Range("H110").Value = 1542.4286
$ws.Range("I110").Value = 1542.4286
$ws.Range("K110").Value = 1542.4286
$ws.Range("M110").Value = 502.5714

$ws.Range("H136").Value = 2217.0715
$ws.Range("I136").Value = 2019.5454
$ws.Range("K136").Value = 6058.6362
$ws.Range("M136").Value = -3508.6362

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 876.4761999999999
$ws.Range("I94").Value = 861.8823
$ws.Range("K94").Value = 861.8823
$ws.Range("M94").Value = -410.8823

$ws.Range("H105").Value = 7167.591
$ws.Range("I105").Value = 11407.909
$ws.Range("K105").Value = 11407.909
$ws.Range("M105").Value = -9660.909

$ws.Range("H128").Value = 4220
$ws.Range("I128").Value = 4220
$ws.Range("K128").Value = 12660
$ws.Range("M128").Value = -10170

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 1889
$ws.Range("J14").Value = 1889
$ws.Range("L14").Value = 1889
$ws.Range("N14").Value = -2229

$ws.Range("H19").Value = 1065
$ws.Range("I19").Value = 385.83334
$ws.Range("J19").Value = 1744.1666
$ws.Range("K19").Value = 385.83334
$ws.Range("L19").Value = 1744.1666
$ws.Range("M19").Value = -215.83334
$ws.Range("N19").Value = -2084.1666

$ws.Range("H24").Value = 1065
$ws.Range("I24").Value = 385.83334
$ws.Range("J24").Value = 1744.1666
$ws.Range("K24").Value = 385.83334
$ws.Range("L24").Value = 1744.1666
$ws.Range("M24").Value = -215.83334
$ws.Range("N24").Value = -2084.1666

$ws.Range("H107").Value = 1554.7097
$ws.Range("I107").Value = 1357.375
$ws.Range("J107").Value = 2231.2856
$ws.Range("K107").Value = 1357.375
$ws.Range("L107").Value = 2231.2856
$ws.Range("M107").Value = 562.625
$ws.Range("N107").Value = -6071.2856

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1096.75
$ws.Range("J107").Value = 1261.1818
$ws.Range("L107").Value = 3783.5454
$ws.Range("N107").Value = -7623.5454

$ws.Range("H131").Value = 1892.05
$ws.Range("I131").Value = 2107.5
$ws.Range("J131").Value = 1838.1875
$ws.Range("K131").Value = 6322.5
$ws.Range("L131").Value = 5514.5625
$ws.Range("M131").Value = -1282.5
$ws.Range("N131").Value = -15594.5625

$ws.Range("H139").Value = 4010
$ws.Range("I139").Value = 4010
$ws.Range("K139").Value = 12030
$ws.Range("M139").Value = -6890

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 28295
$ws.Range("J68").Value = 28295
$ws.Range("L68").Value = 28295
$ws.Range("N68").Value = -29917

$ws.Range("H71").Value = 28295
$ws.Range("J71").Value = 28295
$ws.Range("L71").Value = 84885
$ws.Range("N71").Value = -92997

$ws.Range("H97").Value = 57171.332
$ws.Range("I97").Value = 40751.94
$ws.Range("K97").Value = 40751.94
$ws.Range("M97").Value = -40255.94

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 19049.6
$ws.Range("I7").Value = 23618.334
$ws.Range("J7").Value = 5343.4
$ws.Range("K7").Value = 23618.334
$ws.Range("L7").Value = 5343.4
$ws.Range("M7").Value = -23506.334
$ws.Range("N7").Value = -5567.4

$ws.Range("H25").Value = 4000
$ws.Range("I25").Value = 4000
$ws.Range("K25").Value = 4000
$ws.Range("M25").Value = -3770

$ws.Range("H34").Value = 0
$ws.Range("I34").Value = 0
$ws.Range("K34").Value = 0
$ws.Range("M34").ClearContents()

$ws.Range("H100").Value = 8749.5
$ws.Range("I100").Value = 8749.5
$ws.Range("K100").Value = 8749.5
$ws.Range("M100").Value = -8208.5

$ws.Range("H107").Value = 3253.375
$ws.Range("I107").Value = 3253.375
$ws.Range("K107").Value = 3253.375
$ws.Range("M107").Value = -1333.375

$ws.Range("H126").Value = 19049.6
$ws.Range("I126").Value = 23618.334
$ws.Range("J126").Value = 5343.4
$ws.Range("K126").Value = 70855.00199999999
$ws.Range("L126").Value = 16030.2
$ws.Range("M126").Value = -68385.00199999999
$ws.Range("N126").Value = -20970.2

$ws.Range("H136").Value = 6074.4443
$ws.Range("I136").Value = 4935
$ws.Range("J136").Value = 7498.75
$ws.Range("K136").Value = 14805
$ws.Range("L136").Value = 22496.25
$ws.Range("M136").Value = -12255
$ws.Range("N136").Value = -27596.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 10135899
$ws.Range("I3").Value = 16759865
$ws.Range("K3").Value = 16759865
$ws.Range("M3").Value = -16759751

$ws.Range("H96").Value = 3911.647
$ws.Range("I96").Value = 2499
$ws.Range("K96").Value = 2499
$ws.Range("M96").Value = -1126

$ws.Range("H98").Value = 28590
$ws.Range("J98").Value = 28590
$ws.Range("L98").Value = 28590
$ws.Range("N98").Value = -34580

$ws.Range("H100").Value = 3248.2
$ws.Range("I100").Value = 4079.7
$ws.Range("J100").Value = 2416.7
$ws.Range("K100").Value = 8159.4
$ws.Range("L100").Value = 4833.4
$ws.Range("M100").Value = -7618.4
$ws.Range("N100").Value = -5915.4

$ws.Range("H107").Value = 17858350
$ws.Range("I107").Value = 895.25
$ws.Range("K107").Value = 2685.75
$ws.Range("M107").Value = -765.75

$ws.Range("H136").Value = 3036.9565
$ws.Range("I136").Value = 2720.4546
$ws.Range("K136").Value = 8161.3638
$ws.Range("M136").Value = -5611.3638
